$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for: $old"
    }
}

# 1) "Levantar la hoja que se encuentra encima de las tarjetas." -> "Levantar la hoja."
Replace-Text "Levantar la hoja que se encuentra encima de las tarjetas." "Levantar la hoja."

# 2) "Tomar una de las tarjetas." -> "Poner la hoja en un lado de las tarjetas."
Replace-Text "Tomar una de las tarjetas." "Poner la hoja en un lado de las tarjetas."

# 3) "Tomar la segunda tarjeta." -> "Tomar una de las tarjetas."
Replace-Text "Tomar la segunda tarjeta." "Tomar una de las tarjetas."

# 4) "Poner la hoja nuevamente." -> "Tomar la segunda tarjeta."
Replace-Text "Poner la hoja nuevamente." "Tomar la segunda tarjeta."

# 5) "Ubicar una de las tarjetas sobre la hoja." -> "Poner la hoja nuevamente en la posición original."
Replace-Text "Ubicar una de las tarjetas sobre la hoja." "Poner la hoja nuevamente en la posición original."

# 6) "Poner la otra tarjeta al lado de la otra unidas por el lado más corto." ->
#    "Poner las dos tarjetas de forma perpendicular a la hoja, apoyadas por el lado más corto una tarjeta frente a la otra."
Replace-Text "Poner la otra tarjeta al lado de la otra unidas por el lado más corto." "Poner las dos tarjetas de forma perpendicular a la hoja, apoyadas por el lado más corto una tarjeta frente a la otra."

# 7) "Tomar las dos tarjetas con una mano y ponerlas en forma de triangulo apuntando hacia arriba con los dos lados que se encontraban haciendo contacto." ->
#    "Poner el dedo índice entre las dos tarjetas de forma que exista un espacio entre ellas."
Replace-Text "Tomar las dos tarjetas con una mano y ponerlas en forma de triangulo apuntando hacia arriba con los dos lados que se encontraban haciendo contacto." "Poner el dedo índice entre las dos tarjetas de forma que exista un espacio entre ellas."

# Append 3 new list paragraphs after the last one, inheriting the same list style.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$p8 = $d.Paragraphs($d.Paragraphs.Count)
$p8.Range.Text = "Separar las tarjetas en la parte donde se encuentran apoyadas en la hoja."
$p8.Range.InsertParagraphAfter()

$p9 = $d.Paragraphs($d.Paragraphs.Count)
$p9.Range.Text = "Unir los lados superiores de las tarjetas."
$p9.Range.InsertParagraphAfter()

$p10 = $d.Paragraphs($d.Paragraphs.Count)
$p10.Range.Text = "Retirar la mano, de forma que las tarjetas queden sostenidas por el apoyo que existe entre ellas en la parte superior. (Deben formar un triángulo mirando hacia arriba)"

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
